{"js": "// Update the lattice-multiplication exercise table: each of the 15 cells\n// (5 rows x 3 cols, in reading order) gets a brand-new \"A x B\" problem.\n// Every cell's text is 5 lines joined by manual line breaks:\n//   1) \"A x B\"\n//   2) \"  b1    b2\"   (tens/units digit of B, space-padded)\n//   3) \"  ----\"\n//   4) \"a1|    |\"     (tens digit of A)\n//   5) \"a2|    |\"     (units digit of A)\n\nconst newProblems = [\n  [90, 81], [15, 25], [74, 59],\n  [19, 61], [78, 40], [23, 48],\n  [88, 32], [11, 92], [48, 38],\n  [79, 66], [46, 52], [55, 78],\n  [58, 51], [62, 24], [43, 52],\n];\n\nfunction cellText(a, b) {\n  const aTens = Math.floor(a / 10);\n  const aUnits = a % 10;\n  const bTens = Math.floor(b / 10);\n  const bUnits = b % 10;\n  const lines = [\n    `${a} x ${b}`,\n    `  ${bTens}    ${bUnits}`,\n    `  ----`,\n    `${aTens}|    |`,\n    `${aUnits}|    |`,\n  ];\n  // Office.js represents a manual line break (\"<w:br/>\") inside a\n  // paragraph's text as the vertical-tab character (\\v / \\u000b).\n  return lines.join(\"\\u000b\");\n}\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\nconst rowCount = 5;\nconst colCount = 3;\n\n// Grab every cell's first paragraph up front, then replace its text in a\n// single batch so the existing run formatting (sz=32) is preserved.\nconst paragraphs = [];\nfor (let r = 0; r < rowCount; r++) {\n  for (let c = 0; c < colCount; c++) {\n    const cell = table.getCell(r, c);\n    const paras = cell.body.paragraphs;\n    paras.load(\"items\");\n    paragraphs.push(paras);\n  }\n}\nawait context.sync();\n\nfor (let i = 0; i < newProblems.length; i++) {\n  const [a, b] = newProblems[i];\n  const firstPara = paragraphs[i].items[0];\n  firstPara.insertText(cellText(a, b), Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Update the lattice-multiplication exercise table: each of the 15 cells\n# (5 rows x 3 cols, in reading order) gets a brand-new \"A x B\" problem.\n# Every cell's text is 5 lines joined by manual line breaks:\n#   1) \"A x B\"\n#   2) \"  b1    b2\"   (tens/units digit of B, space-padded)\n#   3) \"  ----\"\n#   4) \"a1|    |\"     (tens digit of A)\n#   5) \"a2|    |\"     (units digit of A)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# New (A, B) pairs, in reading order (row-major), matching the 5x3 table.\n$newProblems = @(\n  @(90, 81), @(15, 25), @(74, 59),\n  @(19, 61), @(78, 40), @(23, 48),\n  @(88, 32), @(11, 92), @(48, 38),\n  @(79, 66), @(46, 52), @(55, 78),\n  @(58, 51), @(62, 24), @(43, 52)\n)\n\n$nl = [char]11  # manual line break (\"<w:br/>\") within a Word paragraph\n\n$rowCount = 5\n$colCount = 3\n\nfor ($i = 0; $i -lt $newProblems.Count; $i++) {\n    $a = $newProblems[$i][0]\n    $b = $newProblems[$i][1]\n\n    $aTens = [math]::Floor($a / 10)\n    $aUnits = $a % 10\n    $bTens = [math]::Floor($b / 10)\n    $bUnits = $b % 10\n\n    $line1 = \"$a x $b\"\n    $line2 = \"  $bTens    $bUnits\"\n    $line3 = \"  ----\"\n    $line4 = \"$aTens|    |\"\n    $line5 = \"$aUnits|    |\"\n\n    $newText = \"$line1$nl$line2$nl$line3$nl$line4$nl$line5\"\n\n    $r = [math]::Floor($i / $colCount) + 1\n    $c = ($i % $colCount) + 1\n\n    $cell = $t.Cell($r, $c)\n    $cell.Range.Text = $newText\n}\n"}
